$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 719.7273
$ws.Range("I41").Value = 92.25
$ws.Range("J41").Value = 1078.2858
$ws.Range("K41").Value = 92.25
$ws.Range("L41").Value = 1078.2858
$ws.Range("M41").Value = 347.75
$ws.Range("N41").Value = -1958.2858
$ws.Range("H86").Value = 7026.1177
$ws.Range("J86").Value = 25704.25
$ws.Range("L86").Value = 25704.25
$ws.Range("N86").Value = -27950.25
$ws.Range("H88").Value = 680.5
$ws.Range("I88").Value = 680.5
$ws.Range("J88").Value = 0
$ws.Range("K88").Value = 680.5
$ws.Range("L88").Value = 0
$ws.Range("M88").Value = -274.5
$ws.Range("N88").ClearContents()
$ws.Range("H89").Value = 7026.1177
$ws.Range("J89").Value = 25704.25
$ws.Range("L89").Value = 128521.25
$ws.Range("N89").Value = -139753.25
$ws.Range("H91").Value = 680.5
$ws.Range("I91").Value = 680.5
$ws.Range("J91").Value = 0
$ws.Range("K91").Value = 680.5
$ws.Range("L91").Value = 0
$ws.Range("M91").Value = 723.5
$ws.Range("N91").ClearContents()
$ws.Range("H98").Value = 917.5
$ws.Range("I98").Value = 957.7143
$ws.Range("J98").Value = 861.2
$ws.Range("K98").Value = 957.7143
$ws.Range("L98").Value = 861.2
$ws.Range("M98").Value = 540.2857
$ws.Range("N98").Value = -3857.2
$ws.Range("H111").Value = 7305.6665
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("N111").ClearContents()
$ws.Range("H112").Value = 1091.1621
$ws.Range("J112").Value = 1107.8
$ws.Range("L112").Value = 3323.4
$ws.Range("N112").Value = -5539.4
$ws.Range("H122").Value = 917.5
$ws.Range("I122").Value = 957.7143
$ws.Range("J122").Value = 861.2
$ws.Range("K122").Value = 2873.1429
$ws.Range("L122").Value = 2583.6
$ws.Range("M122").Value = -423.1428999999998
$ws.Range("N122").Value = -7483.6
$ws.Range("H129").Value = 134271.62
$ws.Range("I129").Value = 575
$ws.Range("J129").Value = 141803.83
$ws.Range("K129").Value = 1725
$ws.Range("L129").Value = 425411.49
$ws.Range("M129").Value = 3275
$ws.Range("N129").Value = -435411.49
$ws.Range("H137").Value = 1384.6666
$ws.Range("I137").Value = 1362.381
$ws.Range("K137").Value = 4087.143
$ws.Range("M137").Value = -1537.143

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").ClearContents()
$ws.Range("H45").Value = 3124.8333
$ws.Range("I45").Value = 3036.5
$ws.Range("J45").Value = 3202.125
$ws.Range("K45").Value = 3036.5
$ws.Range("L45").Value = 3202.125
$ws.Range("M45").Value = -2659.5
$ws.Range("N45").Value = -3956.125

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1762.0571
$ws.Range("I86").Value = 1609.619
$ws.Range("J86").Value = 1990.7142
$ws.Range("K86").Value = 1609.619
$ws.Range("L86").Value = 1990.7142
$ws.Range("M86").Value = -486.6189999999999
$ws.Range("N86").Value = -4236.7142
$ws.Range("H89").Value = 1762.0571
$ws.Range("I89").Value = 1609.619
$ws.Range("J89").Value = 1990.7142
$ws.Range("K89").Value = 8048.094999999999
$ws.Range("L89").Value = 9953.571
$ws.Range("M89").Value = -2432.094999999999
$ws.Range("N89").Value = -21185.571
$ws.Range("H105").Value = 1564254.8
$ws.Range("I105").Value = 1560.7142
$ws.Range("J105").Value = 2779683.2
$ws.Range("K105").Value = 1560.7142
$ws.Range("L105").Value = 2779683.2
$ws.Range("M105").Value = 186.2858000000001
$ws.Range("N105").Value = -2783177.2
$ws.Range("H134").Value = 2987.558
$ws.Range("I134").Value = 2941.9143
$ws.Range("K134").Value = 8825.742899999999
$ws.Range("M134").Value = -6290.742899999999

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 4792.3076
$ws.Range("I62").Value = 4572.727
$ws.Range("K62").Value = 4572.727
$ws.Range("M62").Value = -3948.727
$ws.Range("H65").Value = 4792.3076
$ws.Range("I65").Value = 4572.727
$ws.Range("K65").Value = 22863.635
$ws.Range("M65").Value = -19743.635
$ws.Range("H86").Value = 18041.928
$ws.Range("I86").Value = 11612.667
$ws.Range("J86").Value = 29614.6
$ws.Range("K86").Value = 11612.667
$ws.Range("L86").Value = 29614.6
$ws.Range("M86").Value = -10489.667
$ws.Range("N86").Value = -31860.6
$ws.Range("H89").Value = 18041.928
$ws.Range("I89").Value = 11612.667
$ws.Range("J89").Value = 29614.6
$ws.Range("K89").Value = 58063.335
$ws.Range("L89").Value = 148073
$ws.Range("M89").Value = -52447.335
$ws.Range("N89").Value = -159305
$ws.Range("H137").Value = 24890
$ws.Range("I137").Value = 9000
$ws.Range("K137").Value = 9000
$ws.Range("M137").Value = -3900

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H44").Value = 447.5
$ws.Range("I44").Value = 447.5
$ws.Range("K44").Value = 1342.5
$ws.Range("M44").Value = -944.5
$ws.Range("H50").Value = 199
$ws.Range("I50").Value = 200
$ws.Range("K50").Value = 600
$ws.Range("M50").Value = -119
$ws.Range("H53").Value = 199
$ws.Range("I53").Value = 200
$ws.Range("K53").Value = 600
$ws.Range("M53").Value = -119
$ws.Range("H92").Value = 583.1429000000001
$ws.Range("I92").Value = 260.66666
$ws.Range("J92").Value = 825
$ws.Range("K92").Value = 781.9999799999999
$ws.Range("L92").Value = 2475
$ws.Range("M92").Value = 466.0000200000001
$ws.Range("N92").Value = -4971
$ws.Range("H97").Value = 544.875
$ws.Range("I97").Value = 200
$ws.Range("K97").Value = 600
$ws.Range("M97").Value = -104
$ws.Range("H131").Value = 716.42
$ws.Range("J131").Value = 716.42
$ws.Range("L131").Value = 2149.26
$ws.Range("N131").Value = -12229.26

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3504.8
$ws.Range("I80").Value = 2862.5
$ws.Range("J80").Value = 3933
$ws.Range("K80").Value = 2862.5
$ws.Range("L80").Value = 3933
$ws.Range("M80").Value = -1864.5
$ws.Range("N80").Value = -5929
$ws.Range("H83").Value = 3504.8
$ws.Range("I83").Value = 2862.5
$ws.Range("J83").Value = 3933
$ws.Range("K83").Value = 14312.5
$ws.Range("L83").Value = 19665
$ws.Range("M83").Value = -9320.5
$ws.Range("N83").Value = -29649

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H96").Value = 39500
$ws.Range("J96").Value = 39500
$ws.Range("L96").Value = 39500
$ws.Range("N96").Value = -44992
